$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily")

# New daily measurement log entry for 2022-10-15
$ws.Range("A18").Value = 20221015
$ws.Range("B18").Value = "completed"
$ws.Range("C18").Value = "completed"
$ws.Range("D18").Value = "AH"

# Leave selection where the author left it after logging the entry
$ws.Range("B24").Select()
